$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("eed_t2-dev_t23")
$ws2.Range("B2").Value2 = "z_age2mo_personal_all"
$ws2.Range("D2").Value2 = 139
$ws2.Range("E2").Value2 = 996
$ws2.Range("G2").Value2 = 0.5659999999999999

$ws2.Range("B3").Value2 = "z_age2mo_motor_all"
$ws2.Range("D3").Value2 = 139
$ws2.Range("E3").Value2 = 996
$ws2.Range("G3").Value2 = 0.5659999999999999

$ws2.Range("B4").Value2 = "z_age2mo_combined_all"
$ws2.Range("D4").Value2 = 139
$ws2.Range("E4").Value2 = 996
$ws2.Range("G4").Value2 = 0.5659999999999999

$ws2.Range("B5").Value2 = "z_age2mo_com_all"
$ws2.Range("D5").Value2 = 139
$ws2.Range("E5").Value2 = 996
$ws2.Range("G5").Value2 = 0.5659999999999999

$ws2.Range("B8").Value2 = "z_age2mo_personal_all"
$ws2.Range("D8").Value2 = 139
$ws2.Range("E8").Value2 = 994
$ws2.Range("G8").Value2 = 0.5639999999999999

$ws2.Range("B9").Value2 = "z_age2mo_motor_all"
$ws2.Range("D9").Value2 = 139
$ws2.Range("E9").Value2 = 994
$ws2.Range("G9").Value2 = 0.5639999999999999

$ws2.Range("B10").Value2 = "z_age2mo_combined_all"
$ws2.Range("D10").Value2 = 139
$ws2.Range("E10").Value2 = 994
$ws2.Range("G10").Value2 = 0.5639999999999999

$ws2.Range("B11").Value2 = "z_age2mo_com_all"
$ws2.Range("D11").Value2 = 139
$ws2.Range("E11").Value2 = 994
$ws2.Range("G11").Value2 = 0.5639999999999999

$ws2.Range("B14").Value2 = "z_age2mo_personal_all"
$ws2.Range("D14").Value2 = 139
$ws2.Range("E14").Value2 = 993
$ws2.Range("G14").Value2 = 0.5639999999999999

$ws2.Range("B15").Value2 = "z_age2mo_motor_all"
$ws2.Range("D15").Value2 = 139
$ws2.Range("E15").Value2 = 993
$ws2.Range("G15").Value2 = 0.5639999999999999

$ws2.Range("B16").Value2 = "z_age2mo_combined_all"
$ws2.Range("D16").Value2 = 139
$ws2.Range("E16").Value2 = 993
$ws2.Range("G16").Value2 = 0.5639999999999999

$ws2.Range("B17").Value2 = "z_age2mo_com_all"
$ws2.Range("D17").Value2 = 139
$ws2.Range("E17").Value2 = 993
$ws2.Range("G17").Value2 = 0.5639999999999999

$ws2.Range("B20").Value2 = "z_age2mo_personal_all"
$ws2.Range("D20").Value2 = 139
$ws2.Range("E20").Value2 = 813
$ws2.Range("G20").Value2 = 0.462

$ws2.Range("B21").Value2 = "z_age2mo_motor_all"
$ws2.Range("D21").Value2 = 139
$ws2.Range("E21").Value2 = 813
$ws2.Range("G21").Value2 = 0.462

$ws2.Range("B22").Value2 = "z_age2mo_combined_all"
$ws2.Range("D22").Value2 = 139
$ws2.Range("E22").Value2 = 813
$ws2.Range("G22").Value2 = 0.462

$ws2.Range("B23").Value2 = "z_age2mo_com_all"
$ws2.Range("D23").Value2 = 139
$ws2.Range("E23").Value2 = 813
$ws2.Range("G23").Value2 = 0.462

$ws2.Range("B26").Value2 = "z_age2mo_personal_all"
$ws2.Range("D26").Value2 = 139
$ws2.Range("E26").Value2 = 813
$ws2.Range("G26").Value2 = 0.462

$ws2.Range("B27").Value2 = "z_age2mo_motor_all"
$ws2.Range("D27").Value2 = 139
$ws2.Range("E27").Value2 = 813
$ws2.Range("G27").Value2 = 0.462

$ws2.Range("B28").Value2 = "z_age2mo_combined_all"
$ws2.Range("D28").Value2 = 139
$ws2.Range("E28").Value2 = 813
$ws2.Range("G28").Value2 = 0.462

$ws2.Range("B29").Value2 = "z_age2mo_com_all"
$ws2.Range("D29").Value2 = 139
$ws2.Range("E29").Value2 = 813
$ws2.Range("G29").Value2 = 0.462

$ws2.Range("B32").Value2 = "z_age2mo_personal_all"
$ws2.Range("D32").Value2 = 139
$ws2.Range("E32").Value2 = 420
$ws2.Range("G32").Value2 = 0.239

$ws2.Range("B33").Value2 = "z_age2mo_motor_all"
$ws2.Range("D33").Value2 = 139
$ws2.Range("E33").Value2 = 420
$ws2.Range("G33").Value2 = 0.239

$ws2.Range("B34").Value2 = "z_age2mo_combined_all"
$ws2.Range("D34").Value2 = 139
$ws2.Range("E34").Value2 = 420
$ws2.Range("G34").Value2 = 0.239

$ws2.Range("B35").Value2 = "z_age2mo_com_all"
$ws2.Range("D35").Value2 = 139
$ws2.Range("E35").Value2 = 420
$ws2.Range("G35").Value2 = 0.239

$ws2.Range("B38").Value2 = "z_age2mo_personal_all"
$ws2.Range("D38").Value2 = 139
$ws2.Range("E38").Value2 = 421
$ws2.Range("G38").Value2 = 0.239

$ws2.Range("B39").Value2 = "z_age2mo_motor_all"
$ws2.Range("D39").Value2 = 139
$ws2.Range("E39").Value2 = 421
$ws2.Range("G39").Value2 = 0.239

$ws2.Range("B40").Value2 = "z_age2mo_combined_all"
$ws2.Range("D40").Value2 = 139
$ws2.Range("E40").Value2 = 421
$ws2.Range("G40").Value2 = 0.239

$ws2.Range("B41").Value2 = "z_age2mo_com_all"
$ws2.Range("D41").Value2 = 139
$ws2.Range("E41").Value2 = 421
$ws2.Range("G41").Value2 = 0.239

$ws2.Range("B44").Value2 = "z_age2mo_personal_all"
$ws2.Range("D44").Value2 = 139
$ws2.Range("E44").Value2 = 420
$ws2.Range("G44").Value2 = 0.239

$ws2.Range("B45").Value2 = "z_age2mo_motor_all"
$ws2.Range("D45").Value2 = 139
$ws2.Range("E45").Value2 = 420
$ws2.Range("G45").Value2 = 0.239

$ws2.Range("B46").Value2 = "z_age2mo_combined_all"
$ws2.Range("D46").Value2 = 139
$ws2.Range("E46").Value2 = 420
$ws2.Range("G46").Value2 = 0.239

$ws2.Range("B47").Value2 = "z_age2mo_com_all"
$ws2.Range("D47").Value2 = 139
$ws2.Range("E47").Value2 = 420
$ws2.Range("G47").Value2 = 0.239

$ws2.Range("B50").Value2 = "z_age2mo_personal_all"
$ws2.Range("D50").Value2 = 139
$ws2.Range("E50").Value2 = 456
$ws2.Range("G50").Value2 = 0.259

$ws2.Range("B51").Value2 = "z_age2mo_motor_all"
$ws2.Range("D51").Value2 = 139
$ws2.Range("E51").Value2 = 456
$ws2.Range("G51").Value2 = 0.259

$ws2.Range("B52").Value2 = "z_age2mo_combined_all"
$ws2.Range("D52").Value2 = 139
$ws2.Range("E52").Value2 = 456
$ws2.Range("G52").Value2 = 0.259

$ws2.Range("B53").Value2 = "z_age2mo_com_all"
$ws2.Range("D53").Value2 = 139
$ws2.Range("E53").Value2 = 456
$ws2.Range("G53").Value2 = 0.259

$ws2.Range("B56").Value2 = "z_age2mo_personal_all"
$ws2.Range("D56").Value2 = 139
$ws2.Range("E56").Value2 = 456
$ws2.Range("G56").Value2 = 0.259

$ws2.Range("B57").Value2 = "z_age2mo_motor_all"
$ws2.Range("D57").Value2 = 139
$ws2.Range("E57").Value2 = 456
$ws2.Range("G57").Value2 = 0.259

$ws2.Range("B58").Value2 = "z_age2mo_combined_all"
$ws2.Range("D58").Value2 = 139
$ws2.Range("E58").Value2 = 456
$ws2.Range("G58").Value2 = 0.259

$ws2.Range("B59").Value2 = "z_age2mo_com_all"
$ws2.Range("D59").Value2 = 139
$ws2.Range("E59").Value2 = 456
$ws2.Range("G59").Value2 = 0.259

$ws3 = $wb.Worksheets.Item("reg1b_t2-dev_t23")
$ws3.Range("B2").Value2 = "z_age2mo_personal_all"
$ws3.Range("D2").Value2 = 139
$ws3.Range("E2").Value2 = 422
$ws3.Range("G2").Value2 = 0.24

$ws3.Range("B3").Value2 = "z_age2mo_motor_all"
$ws3.Range("D3").Value2 = 139
$ws3.Range("E3").Value2 = 422
$ws3.Range("G3").Value2 = 0.24

$ws3.Range("B4").Value2 = "z_age2mo_combined_all"
$ws3.Range("D4").Value2 = 139
$ws3.Range("E4").Value2 = 422
$ws3.Range("G4").Value2 = 0.24

$ws3.Range("B5").Value2 = "z_age2mo_com_all"
$ws3.Range("D5").Value2 = 139
$ws3.Range("E5").Value2 = 422
$ws3.Range("G5").Value2 = 0.24
